# "Weitere Ergänzungen für das String-Statement"
# Mark the stringStatement row (row 46) as fully done: set B46 to 1 and
# recolor A46 to match the other "complete" rows (green fill, style index 6),
# then move the view/selection as captured in the saved workbook state.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the completion value for stringStatement (row 46) to 1 (fully done).
$ws.Range("B46").Value = 1

# Re-style A46 the same way the other "done" rows (e.g. A45) are styled -
# copy the format (green fill) from an already-completed row onto A46.
$ws.Range("A45").Copy()
$ws.Range("A46").PasteSpecial(-4122) # xlPasteFormats

# Update the view state: scrolled position and current selection.
$ws.Application.ActiveWindow.ScrollRow = 31
$ws.Range("A45").Select()
